$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "45.576.21"
$ws.Range("E2").Value = "  +7.33%  "

$ws.Range("D3").Value = "2.381.49"
$ws.Range("E3").Value = "  +4.24%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "114.17"
$ws.Range("E5").Value = "  +11.29%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "317.91"
$ws.Range("E6").Value = "  +2.78%  "

$ws.Range("E7").Value = "  +2.12%  "

$ws.Range("E8").Value = "  +0.07%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.627"
$ws.Range("E9").Value = "  +4.55%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.52"
$ws.Range("E10").Value = "  +10.86%  "

$ws.Range("E11").Value = "  +3.66%  "

$ws.Range("E12").Value = "  +6.44%  "

$ws.Range("E13").Value = "  +5.46%  "

$ws.Range("E14").Value = "  +1.74%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.85"
$ws.Range("E15").Value = "  +4.50%  "

$ws.Range("D16").Value = "2.741.81"
$ws.Range("E16").Value = "  +4.23%  "

$ws.Range("D17").Value = "2.380.33"
$ws.Range("E17").Value = "  +3.37%  "

$ws.Range("D18").Value = "45.416.52"
$ws.Range("E18").Value = "  +7.05%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.61"
$ws.Range("E19").Value = "  +4.93%  "

$ws.Range("E20").Value = "  +3.81%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.39"
$ws.Range("E21").Value = "  +0.55%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "74.83"
$ws.Range("E22").Value = "  +2.46%  "

$ws.Range("E23").Value = "  +4.69%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "268.68"
$ws.Range("E24").Value = "  -0.58%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.35"
$ws.Range("E25").Value = "  +9.15%  "

$ws.Range("E26").Value = "  +0.25%  "

$ws.Range("E27").Value = "  +10.57%  "

$ws.Range("E28").Value = "  +6.06%  "

$ws.Range("E29").Value = "  +2.54%  "

$ws.Range("B30").Value = "EthereumClassic"
$ws.Range("C30").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "22.91"
$ws.Range("E30").Value = "  +2.83%  "

$ws.Range("B31").Value = "InjectiveProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "38.76"
$ws.Range("E31").Value = "  +8.91%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0958"
$ws.Range("E32").Value = "  +13.34%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "171.21"
$ws.Range("E33").Value = "  +4.27%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.95"
$ws.Range("E34").Value = "  +16.61%  "

$ws.Range("E35").Value = "  +2.19%  "

$ws.Range("B36").Value = "RenderToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.98"
$ws.Range("E36").Value = "  +11.23%  "

$ws.Range("B37").Value = "Kaspa"
$ws.Range("C37").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.120"
$ws.Range("E37").Value = "  +8.17%  "

$ws.Range("B38").Value = "NEARProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.07"
$ws.Range("E38").Value = "  +13.23%  "

$ws.Range("B39").Value = "LidoDAOToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.04"
$ws.Range("E39").Value = "  +10.73%  "

$ws.Range("E40").Value = "  +6.46%  "

$ws.Range("E41").Value = "  +10.63%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "104.72"
$ws.Range("E42").Value = "  -6.12%  "

$ws.Range("E43").Value = "  +6.89%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "71.40"
$ws.Range("E44").Value = "  +3.18%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.35"
$ws.Range("E45").Value = "  +11.00%  "

$ws.Range("E46").Value = "  +0.15%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "5.82"
$ws.Range("E47").Value = "  +13.11%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "116.35"
$ws.Range("E48").Value = "  +5.96%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.67"
$ws.Range("E49").Value = "  +20.51%  "

$ws.Range("B50").Value = "ordi"
$ws.Range("C50").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "79.77"
$ws.Range("E50").Value = "  +3.64%  "

$ws.Range("B51").Value = "FraxShare"
$ws.Range("C51").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "9.35"
$ws.Range("E51").Value = "  +8.39%  "
